$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 6444.4443
$ws.Range("J40").Value = 9600
$ws.Range("L40").Value = 9600
$ws.Range("N40").Value = -9950

# Row 96
$ws.Range("H96").Value = 358.27274
$ws.Range("I96").Value = 230.25
$ws.Range("J96").Value = 699.6667
$ws.Range("K96").Value = 690.75
$ws.Range("L96").Value = 2099.0001
$ws.Range("M96").Value = 682.25
$ws.Range("N96").Value = -4845.0001

# Row 106
$ws.Range("H106").Value = 2448.5
$ws.Range("I106").Value = 1598
$ws.Range("K106").Value = 1598
$ws.Range("M106").Value = -967

# Row 116
$ws.Range("J116").Value = 5500
$ws.Range("L116").Value = 5500
$ws.Range("N116").Value = -12384

# Row 132
$ws.Range("H132").Value = 3770.4443
$ws.Range("I132").Value = 1622.5834
$ws.Range("J132").Value = 8066.1665
$ws.Range("K132").Value = 4867.7502
$ws.Range("L132").Value = 24198.4995
$ws.Range("M132").Value = -2337.7502
$ws.Range("N132").Value = -29258.4995

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 12020
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 12020
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 12020
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -12594

# Row 88
$ws.Range("H88").Value = 1650
$ws.Range("J88").Value = 1900
$ws.Range("L88").Value = 1900
$ws.Range("N88").Value = -2712

# Row 91
$ws.Range("H91").Value = 1650
$ws.Range("J91").Value = 1900
$ws.Range("L91").Value = 1900
$ws.Range("N91").Value = -4708

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 39
$ws.Range("H39").Value = 28748.25
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 28748.25
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 28748.25
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -29530.25

# Row 49
$ws.Range("H49").Value = 28748.25
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 28748.25
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 28748.25
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -29112.25

# Row 132
$ws.Range("H132").Value = 7428
$ws.Range("I132").Value = 7428
$ws.Range("K132").Value = 22284
$ws.Range("M132").Value = -19754

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

# Row 128
$ws.Range("H128").Value = 98332.664
$ws.Range("I128").Value = 98332.664
$ws.Range("K128").Value = 294997.992
$ws.Range("M128").Value = -290017.992

# Row 131
$ws.Range("H131").Value = 2923.1428
$ws.Range("I131").Value = 4950
$ws.Range("K131").Value = 14850
$ws.Range("M131").Value = -9810

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8990.666999999999
$ws.Range("I70").Value = 8990.666999999999
$ws.Range("K70").Value = 8990.666999999999
$ws.Range("M70").Value = -8720.666999999999

# Row 73
$ws.Range("H73").Value = 8990.666999999999
$ws.Range("I73").Value = 8990.666999999999
$ws.Range("K73").Value = 8990.666999999999
$ws.Range("M73").Value = -8054.666999999999

# Row 97
$ws.Range("H97").Value = 378.63635
$ws.Range("I97").Value = 378.63635
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 378.63635
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 117.36365
$ws.Range("N97").ClearContents()

# Row 99
$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -2754

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3068.5
$ws.Range("I7").Value = 2282.2
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 2282.2
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -2170.2
$ws.Range("N7").Value = -7224

# Row 38
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30820

# Row 82
$ws.Range("H82").Value = 1590.091
$ws.Range("I82").Value = 1874.375
$ws.Range("J82").Value = 832
$ws.Range("K82").Value = 1874.375
$ws.Range("L82").Value = 832
$ws.Range("M82").Value = -1513.375
$ws.Range("N82").Value = -1554

# Row 85
$ws.Range("H85").Value = 1590.091
$ws.Range("I85").Value = 1874.375
$ws.Range("J85").Value = 832
$ws.Range("K85").Value = 1874.375
$ws.Range("L85").Value = 832
$ws.Range("M85").Value = -626.375
$ws.Range("N85").Value = -3328

# Row 100
$ws.Range("H100").Value = 4166.1665
$ws.Range("I100").Value = 4166.1665
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4166.1665
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3625.1665
$ws.Range("N100").ClearContents()

# Row 126
$ws.Range("H126").Value = 3068.5
$ws.Range("I126").Value = 2282.2
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 6846.599999999999
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -4376.599999999999
$ws.Range("N126").Value = -25940

# Row 136
$ws.Range("H136").Value = 4801
$ws.Range("I136").Value = 3604.2
$ws.Range("J136").Value = 6795.6665
$ws.Range("K136").Value = 10812.6
$ws.Range("L136").Value = 20386.9995
$ws.Range("M136").Value = -8262.599999999999
$ws.Range("N136").Value = -25486.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2686.125
$ws.Range("I81").Value = 2098
$ws.Range("J81").Value = 3666.3333
$ws.Range("K81").Value = 4196
$ws.Range("L81").Value = 7332.6666
$ws.Range("M81").Value = -3135
$ws.Range("N81").Value = -9454.6666

# Row 84
$ws.Range("H84").Value = 2686.125
$ws.Range("I84").Value = 2098
$ws.Range("J84").Value = 3666.3333
$ws.Range("K84").Value = 20980
$ws.Range("L84").Value = 36663.333
$ws.Range("M84").Value = -15676
$ws.Range("N84").Value = -47271.333

# Row 113
$ws.Range("H113").Value = 399.66666
$ws.Range("I113").Value = 399.66666
$ws.Range("K113").Value = 1198.99998
$ws.Range("M113").Value = 971.0000199999999

# Row 132
$ws.Range("H132").Value = 4044.5
$ws.Range("I132").Value = 4044.5
$ws.Range("K132").Value = 12133.5
$ws.Range("M132").Value = -9603.5

# Row 136
$ws.Range("H136").Value = 3643.4443
$ws.Range("I136").Value = 3473.875
$ws.Range("K136").Value = 10421.625
$ws.Range("M136").Value = -7871.625
